$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LMS_valid_dataFinalUpdate_data1")

# Update email addresses in column L (rows 2-5)
$ws.Range("L2").Value = "test1@gmail.comcom"
$ws.Range("L3").Value = "test2@gmail.comcom"
$ws.Range("L4").Value = "test3@gmail.comcom"
$ws.Range("L5").Value = "test4@gmail.comcom"

# Update phone numbers in column N (rows 2-5)
$ws.Range("N2").Value = 9999002622
$ws.Range("N3").Value = 6368666733
$ws.Range("N4").Value = 6985749962
$ws.Range("N5").Value = 4866475853

$wb.Save()
